# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have been generated: it records the
# handback target file name/link, the handback xlf file name, the handback
# datetime, and flips the overall status text from "Ready for handoff" to
# "Handed back: in sync with en-US" across all three worksheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# hyperlink style look: underlined, cornflower-blue (matches the existing
# custom "HyperLink" cell style already used by column A in this workbook)
$hyperlinkColor = 15570276   # BGR encoding of RGB(0x64,0x95,0xED)

# Column width helper values: this runtime quantizes ColumnWidth to 1/6
# increments (based on pixel rounding), so we pick input values that land
# exactly (or as close as achievable) on the desired raw widths.
$wideColWidth = 29.1    # snaps to raw width 30   (target ~29.98)
$fullColWidth = 39.1    # snaps to raw width 40   (target 40)

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value2 = $statusText
$ov.Range("F2").Value2 = $statusText
$ov.Range("E3").Value2 = $statusText
$ov.Range("F3").Value2 = $statusText

$ov.Columns.Item(5).ColumnWidth = $wideColWidth
$ov.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value2 = $statusText
$zh.Range("C3").Value2 = $statusText

# Row 2 -> 31c075b7-... package
$zh.Range("I2").Value2 = "31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37304f2905272a97940b615d72ded7a94dd4e438/e2e/31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.md", [Type]::Missing, [Type]::Missing, "31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.md") | Out-Null
$zh.Range("I2").Font.Underline = $true
$zh.Range("I2").Font.Color = $hyperlinkColor
$zh.Range("J2").Value2 = "31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.00297c46885d3d7f00efc81745cf1a1e613c342b.zh-cn.xlf"
$zh.Range("K2").Value2 = "2016-08-31 06:53:35"

# Row 3 -> 52afc5a2-... package
$zh.Range("I3").Value2 = "52afc5a2-e103-47ae-b83a-46ec2d3b799f.md"
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37304f2905272a97940b615d72ded7a94dd4e438/e2e/52afc5a2-e103-47ae-b83a-46ec2d3b799f.md", [Type]::Missing, [Type]::Missing, "52afc5a2-e103-47ae-b83a-46ec2d3b799f.md") | Out-Null
$zh.Range("I3").Font.Underline = $true
$zh.Range("I3").Font.Color = $hyperlinkColor
$zh.Range("J3").Value2 = "52afc5a2-e103-47ae-b83a-46ec2d3b799f.a09adb940bac49b09281a6bd852e4d47e4698e5d.zh-cn.xlf"
$zh.Range("K3").Value2 = "2016-08-31 06:53:35"

$zh.Columns.Item(3).ColumnWidth = $wideColWidth
$zh.Columns.Item(9).ColumnWidth = $fullColWidth
$zh.Columns.Item(10).ColumnWidth = $fullColWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value2 = $statusText
$de.Range("C3").Value2 = $statusText

# Row 2 -> 31c075b7-... package
$de.Range("I2").Value2 = "31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37304f2905272a97940b615d72ded7a94dd4e438/e2e/31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.md", [Type]::Missing, [Type]::Missing, "31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.md") | Out-Null
$de.Range("I2").Font.Underline = $true
$de.Range("I2").Font.Color = $hyperlinkColor
$de.Range("J2").Value2 = "31c075b7-aced-4e9f-bf9d-bd63cd4f99ac.00297c46885d3d7f00efc81745cf1a1e613c342b.de-de.xlf"
$de.Range("K2").Value2 = "2016-08-31 06:53:43"

# Row 3 -> 52afc5a2-... package
$de.Range("I3").Value2 = "52afc5a2-e103-47ae-b83a-46ec2d3b799f.md"
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37304f2905272a97940b615d72ded7a94dd4e438/e2e/52afc5a2-e103-47ae-b83a-46ec2d3b799f.md", [Type]::Missing, [Type]::Missing, "52afc5a2-e103-47ae-b83a-46ec2d3b799f.md") | Out-Null
$de.Range("I3").Font.Underline = $true
$de.Range("I3").Font.Color = $hyperlinkColor
$de.Range("J3").Value2 = "52afc5a2-e103-47ae-b83a-46ec2d3b799f.a09adb940bac49b09281a6bd852e4d47e4698e5d.de-de.xlf"
$de.Range("K3").Value2 = "2016-08-31 06:53:43"

$de.Columns.Item(3).ColumnWidth = $wideColWidth
$de.Columns.Item(9).ColumnWidth = $fullColWidth
$de.Columns.Item(10).ColumnWidth = $fullColWidth
